$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2-25, columns B,D,E,F,G,I,J,K (updated values per commit "case with 380 kV done")
$data = @(
    @(0.02194918437700721, 0.01772501339230459, 0.1048837806690628, 4.3749189305492, 0.002604072062039644, 1.059771133992264, 0.2441256079106608, 2.496919799730904),
    @(0.01911112781201751, 0.01795720512971677, 0.1033732219077628, 4.306034375401339, 0.002610401798504939, 1.049973308013264, 0.2390017085728005, 2.388760583980797),
    @(0.01736058810011087, 0.01811559345419056, 0.1025128017859345, 4.265771460229644, 0.002614488861760647, 1.044213551431518, 0.2359985207071986, 2.324195393637012),
    @(0.01664537101058983, 0.01818413950483055, 0.1021789977891245, 4.249871854495325, 0.002616204997965169, 1.041930248933191, 0.2348104233296411, 2.298345370438199),
    @(0.01652650179131854, 0.01819576388827926, 0.1021245848687826, 4.247262317014304, 0.00261649302437408, 1.041554947851182, 0.2346152923929168, 2.294080736034203),
    @(0.0173509497636104, 0.01811650165861423, 0.1025082319202646, 4.265554980503111, 0.002614511801083726, 1.044182500254585, 0.2359823532302272, 2.323844909424167),
    @(0.02097236013408121, 0.01780180302054113, 0.1043489939713815, 4.350743714262762, 0.002606213037994046, 1.056339335963472, 0.2423290980987076, 2.45924162941094),
    @(0.02800477841608995, 0.017309225730191, 0.1084931448599065, 4.534084476796409, 0.00259152212817315, 1.08223828751067, 0.2559194273350585, 2.739554392142622),
    @(0.03312178879051686, 0.01702195101557535, 0.1118673826937737, 4.67895336033024, 0.002581681690153221, 1.102561316705767, 0.2666179186102084, 2.954774509905462),
    @(0.03543729651386229, 0.01690719625960213, 0.1134748382433983, 4.747119709982115, 0.00257740937253023, 1.112096686044438, 0.2716435259368524, 3.054754942645332),
    @(0.03631221928429795, 0.01686601024783485, 0.1140940288197179, 4.773262289838158, 0.002575820714785919, 1.115749891554003, 0.2735697032102991, 3.092917697288385),
    @(0.0361238763067746, 0.01687477980873275, 0.1139602079287059, 4.767617303482155, 0.00257616156584346, 1.114961214132002, 0.2731538355156573, 3.084685156457738),
    @(0.03550931605992957, 0.01690376250627956, 0.1135255690148753, 4.749263850889264, 0.002577278088914987, 1.112396384775622, 0.2718015295015022, 3.057888533791299),
    @(0.03513262705797615, 0.01692181012633398, 0.1132607069296085, 4.738064850636647, 0.002577965785780249, 1.110830888781202, 0.2709762174480517, 3.041514317539395),
    @(0.03297020575342202, 0.01702976919112942, 0.1117637953534611, 4.674544424636679, 0.002581964989091038, 1.101944052837887, 0.2662927004483748, 2.948282648950226),
    @(0.03164038765765298, 0.01710006284540988, 0.1108640968520369, 4.636159439524477, 0.002584470529282326, 1.096567065259542, 0.2634603672359503, 2.89162238797627),
    @(0.03087436764484863, 0.01714199436302444, 0.1103534343722892, 4.614294351418096, 0.002585930875977208, 1.093501653364022, 0.2618462166282853, 2.85922804487285),
    @(0.03061481452886028, 0.01715645001352861, 0.110181702524752, 4.606927666431773, 0.002586428631469864, 1.09246842565922, 0.2613022500511164, 2.848293276859295),
    @(0.03178206852555832, 0.01709242482986184, 0.1109591650127975, 4.640223520974388, 0.002584201821833842, 1.097136625963678, 0.2637603266707913, 2.897633745646772),
    @(0.03568988031263132, 0.01689518816962732, 0.1136529480553996, 4.754645732741295, 0.002576949348655672, 1.113148582510249, 0.2721981058971323, 3.065751114970055),
    @(0.03823266134382663, 0.01677949916610544, 0.1154746163794158, 4.831349429813741, 0.002572379416967863, 1.123860565903996, 0.2778473828094405, 3.177389631855817),
    @(0.03687660495896239, 0.01684004219404223, 0.1144967453470613, 4.790234074681422, 0.00257480298177479, 1.118120549990081, 0.2748198475640748, 3.117643340803113),
    @(0.03171801923413398, 0.01709587324375406, 0.1109161641871985, 4.638385517889162, 0.002584323242247705, 1.096879046974948, 0.2636246708136554, 2.894915449241125),
    @(0.02611051085585814, 0.0174292836605634, 0.1073144503124581, 4.482717818530261, 0.002595328189133011, 1.075007376761114, 0.2559194273350585, 2.739554392142622)
)

$cols = @("B","D","E","F","G","I","J","K")

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Range($cols[$j] + $row).Value = $data[$i][$j]
    }
}
